$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Available Data" counts for the first four data rows (date, total_cases, new_cases, population)
$ws.Range("B2").Value = 313
$ws.Range("B3").Value = 313
$ws.Range("B4").Value = 313
$ws.Range("B5").Value = 313

# Update rows for Oporavljeni, Testirani, Smrtni sl. with new Available/Missing counts and recompute Missing Pct
$ws.Range("B6").Value = 237
$ws.Range("C6").Value = 76
$ws.Range("D6").Value = 76/237

$ws.Range("B7").Value = 237
$ws.Range("C7").Value = 76
$ws.Range("D7").Value = 76/237

$ws.Range("B8").Value = 237
$ws.Range("C8").Value = 76
$ws.Range("D8").Value = 76/237
